$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two shapes that need their labels reworked (Address -> Task class,
# deletePerson(p) -> deleteTask(t) message) by scanning for their current text,
# since shape index/name alone isn't unique on this slide.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }

    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text

    if ($full -like ":Address*BookParser*") {
        # First paragraph is ":Address" (8 chars) - replace just "Address" so the
        # leading colon run is left untouched, matching how a user would select
        # the word "Address" and retype "Task".
        $c = $tr.Characters(2, 7)
        $c.Text = "Task"
    }
    elseif ($full -eq "deletePerson(p)") {
        # Run 1 is "deletePerson" (12 chars), run 2 is "(p)" (3 chars).
        $run1 = $tr.Characters(1, 12)
        $run1.Text = "deleteTask"

        $run2 = $tr.Characters(11, 3)
        $run2.Text = "(t)"
    }
}
